$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header columns (L1, M1) continuing the existing numeric
# sequence in row 1, copying the formatting (bold + border) of the
# existing header cells so the new cells match the rest of the row.
$ws.Range("K1").Copy()
$ws.Range("L1:M1").PasteSpecial(-4122)

$ws.Range("L1").Value = 10
$ws.Range("M1").Value = 11

# Split the combined "2 F Spanish Teaching" value into two cells and
# shift the remaining row-2 values over to make room, then split the
# combined instructor name "Tungseth-Faber Kim" into two cells as well.
$ws.Range("G2").Value = "Teaching"
$ws.Range("H2").Value = "2 F Spanish"
$ws.Range("I2").Value = "Practicum"
$ws.Range("J2").Value = "TBA"
$ws.Range("K2").Value = "TBA"
$ws.Range("L2").Value = "Tungseth"
$ws.Range("M2").Value = "Faber Kim"
